$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix H39, H40, F41, H41: convert from text "0" to numeric 0 ---
$ws.Range("H39").Value = 0
$ws.Range("H40").Value = 0
$ws.Range("F41").Value = 0
$ws.Range("H41").Value = 0

# --- Append new rows 53-69 (circular dated 2025-10-27) ---

# Row 53
$ws.Range("A53").Value = "'2025-10-27"
$ws.Range("A53").ClearFormats()
$ws.Range("I53").Value = "'207,300"
$ws.Range("I53").ClearFormats()
$ws.Range("D53").Value = "'"
$ws.Range("D53").ClearFormats()
$ws.Range("F53").Value = "'"
$ws.Range("F53").ClearFormats()
$ws.Range("G53").Value = "'"
$ws.Range("G53").ClearFormats()
$ws.Range("H53").Value = "'"
$ws.Range("H53").ClearFormats()
$ws.Range("B53").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice27102025.pdf"
$ws.Range("E53").Value = "327,500 329,000 328,000 327,000 325,500"
$ws.Range("C53").Value = "CHANDERIYA `nLEAD ZINC `nSMELTER"

# Row 54
$ws.Range("A54").Value = "'2025-10-27"
$ws.Range("A54").ClearFormats()
$ws.Range("I54").Value = "'207,300"
$ws.Range("I54").ClearFormats()
$ws.Range("D54").Value = "'"
$ws.Range("D54").ClearFormats()
$ws.Range("F54").Value = "'"
$ws.Range("F54").ClearFormats()
$ws.Range("G54").Value = "'"
$ws.Range("G54").ClearFormats()
$ws.Range("H54").Value = "'"
$ws.Range("H54").ClearFormats()
$ws.Range("B54").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice27102025.pdf"
$ws.Range("C54").Value = "HYDRO-1 UNIT"
$ws.Range("E54").Value = "327,500 329,000 328,000 327,000 325,500"

# Row 55
$ws.Range("A55").Value = "'2025-10-27"
$ws.Range("A55").ClearFormats()
$ws.Range("I55").Value = "'207,300"
$ws.Range("I55").ClearFormats()
$ws.Range("D55").Value = "'"
$ws.Range("D55").ClearFormats()
$ws.Range("F55").Value = "'"
$ws.Range("F55").ClearFormats()
$ws.Range("G55").Value = "'"
$ws.Range("G55").ClearFormats()
$ws.Range("H55").Value = "'"
$ws.Range("H55").ClearFormats()
$ws.Range("B55").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice27102025.pdf"
$ws.Range("E55").Value = "327,500 329,000 328,000 327,000 325,500"
$ws.Range("C55").Value = "NEW HYDRO `nSMELTER `nCHANDERIYA"

# Row 56
$ws.Range("A56").Value = "'2025-10-27"
$ws.Range("A56").ClearFormats()
$ws.Range("D56").Value = "'0"
$ws.Range("D56").ClearFormats()
$ws.Range("E56").Value = "'0"
$ws.Range("E56").ClearFormats()
$ws.Range("H56").Value = "'0"
$ws.Range("H56").ClearFormats()
$ws.Range("I56").Value = "'0"
$ws.Range("I56").ClearFormats()
$ws.Range("F56").Value = "'"
$ws.Range("F56").ClearFormats()
$ws.Range("B56").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice27102025.pdf"
$ws.Range("G56").Value = "0  327,000"
$ws.Range("C56").Value = "ZINC SMELTER `nDEBRI"

# Row 57
$ws.Range("A57").Value = "'2025-10-27"
$ws.Range("A57").ClearFormats()
$ws.Range("H57").Value = "'0"
$ws.Range("H57").ClearFormats()
$ws.Range("I57").Value = "'207,300"
$ws.Range("I57").ClearFormats()
$ws.Range("D57").Value = "'"
$ws.Range("D57").ClearFormats()
$ws.Range("F57").Value = "'"
$ws.Range("F57").ClearFormats()
$ws.Range("G57").Value = "'"
$ws.Range("G57").ClearFormats()
$ws.Range("B57").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice27102025.pdf"
$ws.Range("E57").Value = "327,500 329,000 328,000 327,000"
$ws.Range("C57").Value = "Pantnagar `nMelting&Castin `ngPlant"

# Row 58
$ws.Range("A58").Value = "'2025-10-27"
$ws.Range("A58").ClearFormats()
$ws.Range("D58").Value = "'0"
$ws.Range("D58").ClearFormats()
$ws.Range("E58").Value = "'0"
$ws.Range("E58").ClearFormats()
$ws.Range("F58").Value = "'0"
$ws.Range("F58").ClearFormats()
$ws.Range("G58").Value = "'0"
$ws.Range("G58").ClearFormats()
$ws.Range("H58").Value = "'0"
$ws.Range("H58").ClearFormats()
$ws.Range("I58").Value = "'207,300"
$ws.Range("I58").ClearFormats()
$ws.Range("B58").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice27102025.pdf"
$ws.Range("C58").Value = "RAJPURA DARIBA `nLEAD SMELTER"

# Row 59
$ws.Range("A59").Value = "'2025-10-27"
$ws.Range("A59").ClearFormats()
$ws.Range("I59").Value = "'209,800"
$ws.Range("I59").ClearFormats()
$ws.Range("D59").Value = "'"
$ws.Range("D59").ClearFormats()
$ws.Range("F59").Value = "'"
$ws.Range("F59").ClearFormats()
$ws.Range("G59").Value = "'"
$ws.Range("G59").ClearFormats()
$ws.Range("H59").Value = "'"
$ws.Range("H59").ClearFormats()
$ws.Range("B59").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice27102025.pdf"
$ws.Range("E59").Value = "330,000 331,500 325,500 329,500 328,000"
$ws.Range("C59").Value = "Faridabad `nDepot"

# Row 60
$ws.Range("A60").Value = "'2025-10-27"
$ws.Range("A60").ClearFormats()
$ws.Range("I60").Value = "'210,200"
$ws.Range("I60").ClearFormats()
$ws.Range("D60").Value = "'"
$ws.Range("D60").ClearFormats()
$ws.Range("F60").Value = "'"
$ws.Range("F60").ClearFormats()
$ws.Range("G60").Value = "'"
$ws.Range("G60").ClearFormats()
$ws.Range("H60").Value = "'"
$ws.Range("H60").ClearFormats()
$ws.Range("B60").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice27102025.pdf"
$ws.Range("C60").Value = "Panvel Depot"
$ws.Range("E60").Value = "330,800 332,300 331,300 330,300 328,800"

# Row 61
$ws.Range("A61").Value = "'2025-10-27"
$ws.Range("A61").ClearFormats()
$ws.Range("I61").Value = "'210,600"
$ws.Range("I61").ClearFormats()
$ws.Range("D61").Value = "'"
$ws.Range("D61").ClearFormats()
$ws.Range("F61").Value = "'"
$ws.Range("F61").ClearFormats()
$ws.Range("G61").Value = "'"
$ws.Range("G61").ClearFormats()
$ws.Range("H61").Value = "'"
$ws.Range("H61").ClearFormats()
$ws.Range("B61").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice27102025.pdf"
$ws.Range("C61").Value = "Pune Depot"
$ws.Range("E61").Value = "330,800 332,300 331,300 330,300 328,800"

# Row 62
$ws.Range("A62").Value = "'2025-10-27"
$ws.Range("A62").ClearFormats()
$ws.Range("I62").Value = "'210,600"
$ws.Range("I62").ClearFormats()
$ws.Range("D62").Value = "'"
$ws.Range("D62").ClearFormats()
$ws.Range("F62").Value = "'"
$ws.Range("F62").ClearFormats()
$ws.Range("G62").Value = "'"
$ws.Range("G62").ClearFormats()
$ws.Range("H62").Value = "'"
$ws.Range("H62").ClearFormats()
$ws.Range("B62").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice27102025.pdf"
$ws.Range("C62").Value = "Baroda Depot"
$ws.Range("E62").Value = "330,800 332,300 331,300 330,300 328,800"

# Row 63
$ws.Range("A63").Value = "'2025-10-27"
$ws.Range("A63").ClearFormats()
$ws.Range("I63").Value = "'210,600"
$ws.Range("I63").ClearFormats()
$ws.Range("D63").Value = "'"
$ws.Range("D63").ClearFormats()
$ws.Range("F63").Value = "'"
$ws.Range("F63").ClearFormats()
$ws.Range("G63").Value = "'"
$ws.Range("G63").ClearFormats()
$ws.Range("H63").Value = "'"
$ws.Range("H63").ClearFormats()
$ws.Range("B63").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice27102025.pdf"
$ws.Range("C63").Value = "Raipur Depot"
$ws.Range("E63").Value = "330,800 332,300 331,300 330,300 328,800"

# Row 64
$ws.Range("A64").Value = "'2025-10-27"
$ws.Range("A64").ClearFormats()
$ws.Range("I64").Value = "'208,300"
$ws.Range("I64").ClearFormats()
$ws.Range("D64").Value = "'"
$ws.Range("D64").ClearFormats()
$ws.Range("F64").Value = "'"
$ws.Range("F64").ClearFormats()
$ws.Range("G64").Value = "'"
$ws.Range("G64").ClearFormats()
$ws.Range("H64").Value = "'"
$ws.Range("H64").ClearFormats()
$ws.Range("B64").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice27102025.pdf"
$ws.Range("E64").Value = "328,500 330,000 329,000 328,000 326,500"
$ws.Range("C64").Value = "JAMSHEDPUR `nSTOCK POINT"

# Row 65
$ws.Range("A65").Value = "'2025-10-27"
$ws.Range("A65").ClearFormats()
$ws.Range("I65").Value = "'208,300"
$ws.Range("I65").ClearFormats()
$ws.Range("C65").Value = "'"
$ws.Range("C65").ClearFormats()
$ws.Range("E65").Value = "'"
$ws.Range("E65").ClearFormats()
$ws.Range("F65").Value = "'"
$ws.Range("F65").ClearFormats()
$ws.Range("G65").Value = "'"
$ws.Range("G65").ClearFormats()
$ws.Range("H65").Value = "'"
$ws.Range("H65").ClearFormats()
$ws.Range("B65").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice27102025.pdf"
$ws.Range("D65").Value = "Kolkata Depot  328,500 330,000 329,000 328,000 326,500"

# Row 66
$ws.Range("A66").Value = "'2025-10-27"
$ws.Range("A66").ClearFormats()
$ws.Range("I66").Value = "'208,300"
$ws.Range("I66").ClearFormats()
$ws.Range("D66").Value = "'"
$ws.Range("D66").ClearFormats()
$ws.Range("F66").Value = "'"
$ws.Range("F66").ClearFormats()
$ws.Range("G66").Value = "'"
$ws.Range("G66").ClearFormats()
$ws.Range("H66").Value = "'"
$ws.Range("H66").ClearFormats()
$ws.Range("B66").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice27102025.pdf"
$ws.Range("E66").Value = "328,500 330,000 329,000 328,000 326,500"
$ws.Range("C66").Value = "Bangalore `nDepot"

# Row 67
$ws.Range("A67").Value = "'2025-10-27"
$ws.Range("A67").ClearFormats()
$ws.Range("I67").Value = "'208,300"
$ws.Range("I67").ClearFormats()
$ws.Range("D67").Value = "'"
$ws.Range("D67").ClearFormats()
$ws.Range("F67").Value = "'"
$ws.Range("F67").ClearFormats()
$ws.Range("G67").Value = "'"
$ws.Range("G67").ClearFormats()
$ws.Range("H67").Value = "'"
$ws.Range("H67").ClearFormats()
$ws.Range("B67").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice27102025.pdf"
$ws.Range("E67").Value = "328,500 330,000 329,000 328,000 326,500"
$ws.Range("C67").Value = "Hyderabad `nDepot"

# Row 68
$ws.Range("A68").Value = "'2025-10-27"
$ws.Range("A68").ClearFormats()
$ws.Range("I68").Value = "'208,300"
$ws.Range("I68").ClearFormats()
$ws.Range("C68").Value = "'"
$ws.Range("C68").ClearFormats()
$ws.Range("E68").Value = "'"
$ws.Range("E68").ClearFormats()
$ws.Range("F68").Value = "'"
$ws.Range("F68").ClearFormats()
$ws.Range("G68").Value = "'"
$ws.Range("G68").ClearFormats()
$ws.Range("H68").Value = "'"
$ws.Range("H68").ClearFormats()
$ws.Range("B68").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice27102025.pdf"
$ws.Range("D68").Value = "Chennai Depot  328,500 330,000 329,000 328,000 326,500"

# Row 69
$ws.Range("A69").Value = "'2025-10-27"
$ws.Range("A69").ClearFormats()
$ws.Range("I69").Value = "'207,300"
$ws.Range("I69").ClearFormats()
$ws.Range("D69").Value = "'"
$ws.Range("D69").ClearFormats()
$ws.Range("F69").Value = "'"
$ws.Range("F69").ClearFormats()
$ws.Range("H69").Value = "'"
$ws.Range("H69").ClearFormats()
$ws.Range("B69").Value = "https://rashtriyametal.com/wp-content/uploads/2025/10/ListPrice27102025.pdf"
$ws.Range("E69").Value = "327,500 329,000"
$ws.Range("G69").Value = "0  327,000 325,500"
$ws.Range("C69").Value = "Sindesar `nsmelter HZAPL"
